# This script re-applies the latest scheduled market-data refresh to the
# Marilith Profits leve-crafting workbook. Columns H-N on the affected rows
# hold live Market Board pricing pulled in by the scheduled runner
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]); this
# script writes the refreshed figures back into each crafting job's sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2283.3333
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2283.3333
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6849.999899999999
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -7185.999899999999
$ws.Range("H46").Value = 7344.6665
$ws.Range("I46").Value = 7344.6665
$ws.Range("K46").Value = 22033.9995
$ws.Range("M46").Value = -21914.9995
$ws.Range("H58").Value = 1313.5
$ws.Range("I58").Value = 1313.5
$ws.Range("K58").Value = 3940.5
$ws.Range("M58").Value = -3790.5
$ws.Range("H60").Value = 7344.6665
$ws.Range("I60").Value = 7344.6665
$ws.Range("K60").Value = 22033.9995
$ws.Range("M60").Value = -21549.9995
$ws.Range("H64").Value = 3983.5
$ws.Range("I64").Value = 3550
$ws.Range("K64").Value = 3550
$ws.Range("M64").Value = -3302
$ws.Range("H67").Value = 3983.5
$ws.Range("I67").Value = 3550
$ws.Range("K67").Value = 3550
$ws.Range("M67").Value = -2692
$ws.Range("H86").Value = 102999.86
$ws.Range("I86").Value = 90499.5
$ws.Range("J86").Value = 108000
$ws.Range("K86").Value = 90499.5
$ws.Range("L86").Value = 108000
$ws.Range("M86").Value = -89376.5
$ws.Range("N86").Value = -110246
$ws.Range("H87").Value = 39995
$ws.Range("J87").Value = 39995
$ws.Range("L87").Value = 39995
$ws.Range("N87").Value = -42491
$ws.Range("H89").Value = 102999.86
$ws.Range("I89").Value = 90499.5
$ws.Range("J89").Value = 108000
$ws.Range("K89").Value = 452497.5
$ws.Range("L89").Value = 540000
$ws.Range("M89").Value = -446881.5
$ws.Range("N89").Value = -551232
$ws.Range("H90").Value = 39995
$ws.Range("J90").Value = 39995
$ws.Range("L90").Value = 119985
$ws.Range("N90").Value = -132465
$ws.Range("H92").Value = 360.375
$ws.Range("I92").Value = 319
$ws.Range("K92").Value = 319
$ws.Range("M92").Value = 929
$ws.Range("H113").Value = 7021.6924
$ws.Range("J113").Value = 7826.143
$ws.Range("L113").Value = 7826.143
$ws.Range("N113").Value = -14334.143
$ws.Range("H137").Value = 3611.4443
$ws.Range("I137").Value = 3611.4443
$ws.Range("K137").Value = 10834.3329
$ws.Range("M137").Value = -8284.332900000001
$ws.Range("H138").Value = 3507.7368
$ws.Range("I138").Value = 1911.75
$ws.Range("J138").Value = 3933.3333
$ws.Range("K138").Value = 5735.25
$ws.Range("L138").Value = 11799.9999
$ws.Range("M138").Value = -595.25
$ws.Range("N138").Value = -22079.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1769.738
$ws.Range("I32").Value = 1358.25
$ws.Range("J32").Value = 9999.5
$ws.Range("K32").Value = 1358.25
$ws.Range("L32").Value = 9999.5
$ws.Range("M32").Value = -1071.25
$ws.Range("N32").Value = -10573.5
$ws.Range("H97").Value = 739.8333
$ws.Range("I97").Value = 739.8333
$ws.Range("K97").Value = 739.8333
$ws.Range("M97").Value = -243.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 12038.4
$ws.Range("I16").Value = 2730.6667
$ws.Range("J16").Value = 26000
$ws.Range("K16").Value = 2730.6667
$ws.Range("L16").Value = 26000
$ws.Range("M16").Value = -2560.6667
$ws.Range("N16").Value = -26340
$ws.Range("H20").Value = 8504
$ws.Range("I20").Value = 8504
$ws.Range("K20").Value = 8504
$ws.Range("M20").Value = -8257
$ws.Range("H86").Value = 2800.5
$ws.Range("I86").Value = 2989.4443
$ws.Range("K86").Value = 2989.4443
$ws.Range("M86").Value = -1866.4443
$ws.Range("H89").Value = 2800.5
$ws.Range("I89").Value = 2989.4443
$ws.Range("K89").Value = 14947.2215
$ws.Range("M89").Value = -9331.2215
$ws.Range("H99").Value = 4599.8
$ws.Range("J99").Value = 4999.5
$ws.Range("L99").Value = 4999.5
$ws.Range("N99").Value = -7995.5
$ws.Range("H105").Value = 1701.6428
$ws.Range("I105").Value = 1591.4445
$ws.Range("K105").Value = 1591.4445
$ws.Range("M105").Value = 155.5554999999999
$ws.Range("H134").Value = 4768.543
$ws.Range("I134").Value = 5455.7144
$ws.Range("J134").Value = 3737.7856
$ws.Range("K134").Value = 16367.1432
$ws.Range("L134").Value = 11213.3568
$ws.Range("M134").Value = -13832.1432
$ws.Range("N134").Value = -16283.3568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 10000
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("M69").Value = -9251
$ws.Range("N69").Value = -11498
$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 10000
$ws.Range("K72").Value = 30000
$ws.Range("L72").Value = 30000
$ws.Range("M72").Value = -26256
$ws.Range("N72").Value = -37488
$ws.Range("H107").Value = 567.5714
$ws.Range("I107").Value = 588
$ws.Range("J107").Value = 516.5
$ws.Range("K107").Value = 588
$ws.Range("L107").Value = 516.5
$ws.Range("M107").Value = 1332
$ws.Range("N107").Value = -4356.5
$ws.Range("H132").Value = 5062.773
$ws.Range("I132").Value = 5600
$ws.Range("J132").Value = 5009.05
$ws.Range("K132").Value = 16800
$ws.Range("L132").Value = 15027.15
$ws.Range("M132").Value = -14270
$ws.Range("N132").Value = -20087.15

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H128").Value = 308312
$ws.Range("I128").Value = 308312
$ws.Range("K128").Value = 924936
$ws.Range("M128").Value = -919956
$ws.Range("H133").Value = 950
$ws.Range("I133").Value = 950
$ws.Range("K133").Value = 2850
$ws.Range("M133").Value = 2210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4003.2
$ws.Range("J22").Value = 6002.6665
$ws.Range("L22").Value = 6002.6665
$ws.Range("N22").Value = -7060.6665
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H126").Value = 10730.143
$ws.Range("I126").Value = 5932.3335
$ws.Range("J126").Value = 14328.5
$ws.Range("K126").Value = 17797.0005
$ws.Range("L126").Value = 42985.5
$ws.Range("M126").Value = -15327.0005
$ws.Range("N126").Value = -47925.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 330.54544
$ws.Range("J55").Value = 362.6
$ws.Range("L55").Value = 362.6
$ws.Range("N55").Value = -708.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 22200
$ws.Range("J101").Value = 22200
$ws.Range("L101").Value = 22200
$ws.Range("N101").Value = -28690
$ws.Range("H104").Value = 38666.332
$ws.Range("J104").Value = 38666.332
$ws.Range("L104").Value = 38666.332
$ws.Range("N104").Value = -45654.332
$ws.Range("H113").Value = 413
$ws.Range("I113").Value = 391.25
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 1173.75
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 996.25
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 3583.1428
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 4047
$ws.Range("K132").Value = 2400
$ws.Range("L132").Value = 12141
$ws.Range("M132").Value = 130
$ws.Range("N132").Value = -17201
$ws.Range("H136").Value = 8282.23
$ws.Range("I136").Value = 6256.857
$ws.Range("J136").Value = 10645.167
$ws.Range("K136").Value = 18770.571
$ws.Range("L136").Value = 31935.501
$ws.Range("M136").Value = -16220.571
$ws.Range("N136").Value = -37035.501
